# [Jallal] more code clean up
# Expand the aggregated "Carrier Route" counts into the more granular
# per-route breakdown, inserting new routes and updating counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New route/count rows (A = Carrier Route code, B = Count), in final order.
$routes = @(
    @("48108-C024", 1),
    @("48105-R007", 1),
    @("48108-R029", 1),
    @("48108-R005", 1),
    @("48108-C090", 4),
    @("48108-R015", 2),
    @("48103-C030", 1)
)

$row = 2
foreach ($entry in $routes) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row = $row + 1
}
